# Update the two-digit multiplication problems in the document.
# Each "old" equation text is unique within the document, so a simple
# Find/Replace (replacing only the first/only match) for each pair is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("96×41=", "70×67="),
    @("53×54=", "48×65="),
    @("37×71=", "35×90="),
    @("98×87=", "69×81="),
    @("46×33=", "48×50="),
    @("24×66=", "81×56="),
    @("62×77=", "37×70="),
    @("51×52=", "33×22="),
    @("58×85=", "76×17="),
    @("63×92=", "18×99="),
    @("94×95=", "90×74="),
    @("92×66=", "39×69="),
    @("27×32=", "40×78="),
    @("30×84=", "56×23="),
    @("87×94=", "91×26="),
    @("25×95=", "53×87="),
    @("82×74=", "43×14="),
    @("84×71=", "47×35="),
    @("37×17=", "37×70="),
    @("41×73=", "28×25="),
    @("31×35=", "22×77="),
    @("99×90=", "28×99="),
    @("11×62=", "40×16="),
    @("19×22=", "53×69="),
    @("87×91=", "99×59=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
